$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.630.51"
$ws.Range("E2").Value = "  +0.04%  "
$ws.Range("D3").Value = "2.298.48"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.34%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.62"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.19%  "
$ws.Range("E7").Value = "  +0.61%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.492"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.54"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.15"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +5.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0782"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.72"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.72%  "
$ws.Range("D15").Value = "2.653.74"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("D16").Value = "2.299.05"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.779"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "42.572.23"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.14%  "
$ws.Range("E20").Value = "  -0.12%  "
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.63"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +7.87%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.49"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("E28").Value = "  +14.95%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "164.92"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.52%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "31.77"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.88%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.65"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.86%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.41"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.22%  "
$ws.Range("E36").Value = "  -1.96%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("E39").Value = "  +0.02%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.08"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +12.80%  "
$ws.Range("D43").Value = "1.948.49"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.08"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.50%  "
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("D48").Value = "2.522.40"
$ws.Range("E48").Value = "  +0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.07"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  -3.95%  "
$ws.Range("E51").Value = "  +1.28%  "
